# Validation_Template_v1.0_.xlsx edit:
#   "added fail/exterior logic removed incomplete scenes"
#
# Adds a new validation row describing a catch-all "All" store type that
# only requires a single Exterior Store Photo image (0 scenes, 1 image,
# target 1), using a wrapped-text style that matches the rest of the
# "template" column but with the sheet's default (non-red) font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 6)
$ws.Range("A6").Value = "All"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "Exterior Store Photo"
$ws.Range("E6").Value = 1

# Match the look of the other "template" cells (D2:D5): wrap text so long
# descriptions remain readable, but keep the default worksheet font rather
# than the red Verdana font used for the free-form template notes.
$ws.Range("D6").Font.Name = "Arial"
$ws.Range("D6").Font.Size = 10
$ws.Range("D6").WrapText = $true

# Leave the cursor where the author's session ended.
$ws.Range("AI47").Select() | Out-Null
